# Update scraped convention/event statistics ("output generated at 456a3b4").
# Sheet "展览" (Exhibitions): row 2 "想去人数" 4338 -> 4374, row 7 "想去人数" 44 -> 45
# Sheet "全部类型" (All types): row 2 "想去人数" 4338 -> 4374, row 8 "想去人数" 44 -> 45

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item("展览")
$wsExhibitions.Range("F2").Value = 4374
$wsExhibitions.Range("F7").Value = 45

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 4374
$wsAllTypes.Range("F8").Value = 45
